$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(14, 20, 23, 35, 37)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = -1
}
